$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE3").Value = 6500
$ws.Range("AE4").Value = 3500
$ws.Range("AE5").Value = 5000
$ws.Range("AE6").Value = 400
$ws.Range("AE7").Value = 98894
$ws.Range("AE8").Value = 281859.52
$ws.Range("AE9").Value = 345
$ws.Range("AE10").Value = 277842
$ws.Range("AE12").Value = 9235
$ws.Range("AE13").Value = 2024
$ws.Range("AE15").Value = 100
$ws.Range("AE18").Value = 150
$ws.Range("AE19").Value = 218
$ws.Range("AE20").Value = 276
$ws.Range("AE21").Value = 60
$ws.Range("AE23").Value = 1190
$ws.Range("AE24").Value = 2835.54
$ws.Range("AE25").Value = 39850
$ws.Range("AE27").Value = 2651
$ws.Range("AE30").Value = 1700
$ws.Range("AE31").Value = 14500
$ws.Range("AE32").Value = 1624.78
$ws.Range("AE33").Value = 1016.52
$ws.Range("AE34").Value = 13000
$ws.Range("AE35").Value = 2000
$ws.Range("AE36").Value = 10000
$ws.Range("AE37").Value = 30000
$ws.Range("AE38").Value = 806
$ws.Range("AE39").Value = 3200
$ws.Range("AE40").Value = 1600
$ws.Range("AE41").Value = 7500
$ws.Range("AE42").Value = 7500
$ws.Range("AE43").Value = 25300
$ws.Range("AE44").Value = 22000
$ws.Range("AE45").Value = 10
$ws.Range("AE46").Value = 1350
$ws.Range("AE47").Value = 1122
$ws.Range("AE48").Value = 1767.64
$ws.Range("AE49").Value = 50272
$ws.Range("AE50").Value = 48062
$ws.Range("AE51").Value = 54691
$ws.Range("AE52").Value = 50272
$ws.Range("AE53").Value = 48062
$ws.Range("AE54").Value = 45133
$ws.Range("AE55").Value = 49868
$ws.Range("AE56").Value = 51762
$ws.Range("AE57").Value = 46080
$ws.Range("AE58").Value = 41346
$ws.Range("AE59").Value = 40399
$ws.Range("AE60").Value = 44187
$ws.Range("AE61").Value = 44187
$ws.Range("AE62").Value = 42293
$ws.Range("AE63").Value = 46500
$ws.Range("AE64").Value = 20000
$ws.Range("AE65").Value = 2839
$ws.Range("AE66").Value = 5037
$ws.Range("AE67").Value = 6815
$ws.Range("AE68").Value = 5037
$ws.Range("AE69").Value = 1500
$ws.Range("AE70").Value = 1000
$ws.Range("AE71").Value = 1000
$ws.Range("AE72").Value = 1000
$ws.Range("AE73").Value = 1500
$ws.Range("AE74").Value = 7851
$ws.Range("AE75").Value = 1989
$ws.Range("AE76").Value = 5882
$ws.Range("AE77").Value = 427
$ws.Range("AE78").Value = 33726
$ws.Range("AE79").Value = 18042
$ws.Range("AE80").Value = 18042
$ws.Range("AE81").Value = 39832
$ws.Range("AE82").Value = 28367
$ws.Range("AE83").Value = 28980
$ws.Range("AE84").Value = 31255
$ws.Range("AE85").Value = 31255
$ws.Range("AE86").Value = 26442
$ws.Range("AE87").Value = 32042
$ws.Range("AE88").Value = 36155
$ws.Range("AE89").Value = 23642
$ws.Range("AE90").Value = 27846
$ws.Range("AE91").Value = 30555
$ws.Range("AE92").Value = 26267
$ws.Range("AE93").Value = 6927
$ws.Range("AE94").Value = 6927
$ws.Range("AE95").Value = 9371
$ws.Range("AE96").Value = 1487
$ws.Range("AE97").Value = 1500
$ws.Range("AE98").Value = 850
$ws.Range("AE100").Value = 5200
$ws.Range("AE101").Value = 7550
$ws.Range("AE102").Value = 5546
$ws.Range("AE103").Value = 5546
$ws.Range("AE105").Value = 27796.23
$ws.Range("AE107").Value = 3354

$ws.Rows.Item(108).Delete()
